$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.705.33"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").Value = "3.556.68"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.75"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.86"
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +2.14%  "
$ws.Range("D8").Value = "3.547.10"
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("E10").Value = "  +18.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.72"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("E13").Value = "  +5.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.50"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").Value = "4.124.67"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "70.765.73"
$ws.Range("E16").Value = "  +2.70%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.15"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.76"
$ws.Range("E18").Value = "  +4.01%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.553.54"
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "571.95"
$ws.Range("E20").Value = "  +5.54%  "
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.72"
$ws.Range("E23").Value = "  -5.50%  "
$ws.Range("E24").Value = "  +3.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.90"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "93.70"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.25"
$ws.Range("E27").Value = "  +4.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.96"
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.26"
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.46"
$ws.Range("E30").Value = "  +2.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.21"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.31"
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.06"
$ws.Range("E34").Value = "  -2.72%  "
$ws.Range("E35").Value = "  +11.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.68"
$ws.Range("E36").Value = "  +16.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "542.81"
$ws.Range("E37").Value = "  -4.01%  "
$ws.Range("E38").Value = "  +3.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.22"
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").Value = "0.0₃0804"
$ws.Range("E40").Value = "  +4.89%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "3.575.64"
$ws.Range("E42").Value = "  +10.70%  "
$ws.Range("E43").Value = "  +3.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.44"
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0469"
$ws.Range("E45").Value = "  +6.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.49"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("E48").Value = "  +3.89%  "
$ws.Range("E49").Value = "  +2.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.49"
$ws.Range("E50").Value = "  +10.61%  "
$ws.Range("E51").Value = "  +0.09%  "
